$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is Saicharan Gurudu (s531499@nwmissouri.edu) - this student record is
# removed entirely. Deleting the whole row shifts every row below it up by one,
# so Sravya Kancharla's row (previously row 9) now lands on row 8.
$ws.Rows(8).Delete()

# Sravya Kancharla's old email (s531500@nwmissouri.edu) is being replaced with a
# hinted/obfuscated address. Move her record down to the bottom of the list
# (row 21): delete her current row (8), which shifts rows 9-21 up to 8-20, then
# re-create her record as the new last row with the updated, hyperlinked email.
$ws.Rows(8).Delete()
$ws.Range("A21").Value = "s@nwmissouri.edu"
$ws.Range("B21").Value = "Sravya Kancharla"
$ws.Hyperlinks.Add($ws.Range("A21"), "mailto:s@nwmissouri.edu")

# Match the author's final view/selection state.
$excel.ActiveWindow.ScrollRow = 4
$null = $ws.Range("B23").Select()
